$wb = $excel.ActiveWorkbook

# --- Insert a new "Login_Valid" sheet before the first sheet ("Login") ---
$loginSheet = $wb.Worksheets.Item("Login")
$newSheet = $wb.Worksheets.Add($loginSheet)
$newSheet.Name = "Login_Valid"

# Header row
$newSheet.Range("A1").Value = "username"
$newSheet.Range("B1").Value = "password"

# Data row (re-uses existing Codecrackers / Numpyninja@2025 values that used
# to live on the Login sheet)
$newSheet.Range("A2").Value = "Codecrackers"
$newSheet.Range("B2").Value = "Numpyninja@2025"
$newSheet.Range("B2").Style = "Hyperlink"
$newSheet.Hyperlinks.Add($newSheet.Range("B2"), "mailto:Numpyninja@2025") | Out-Null

# --- Clear the old data on the Login sheet that moved to Login_Valid ---
$loginSheet = $wb.Worksheets.Item("Login")
$loginSheet.Range("A6").ClearContents()
$loginSheet.Range("B6").ClearContents()

# --- Make the Login sheet the active tab, with the new selection ---
$loginSheet.Activate()
$loginSheet.Range("A12").Select()
